$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.663.34'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '2.246.19'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = "'258.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.49%  '
$ws.Range('D6').Value = "'79.24"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.49%  '
$ws.Range('D7').Value = "'0.624"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.34%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('D9').Value = "'0.605"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('D10').Value = "'43.50"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.76%  '
$ws.Range('D11').Value = "'0.0927"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('D12').Value = "'7.12"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.66%  '
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('D14').Value = '2.567.72'
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').Value = "'14.70"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').Value = '2.246.48'
$ws.Range('E16').Value = '  +0.72%  '
$ws.Range('D17').Value = "'0.798"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('D18').Value = '43.566.77'
$ws.Range('E18').Value = '  +1.62%  '
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').Value = "'71.54"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.50%  '
$ws.Range('D21').Value = "'6.06"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.03%  '
$ws.Range('D22').Value = "'2.33"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.90%  '
$ws.Range('D23').Value = "'232.43"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').Value = "'9.39"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = "'42.24"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.16%  '
$ws.Range('D27').Value = "'10.88"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.76%  '
$ws.Range('E28').Value = '  -2.10%  '
$ws.Range('D29').Value = "'2.23"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('D31').Value = "'173.48"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.93%  '
$ws.Range('D32').Value = "'20.57"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.96%  '
$ws.Range('D33').Value = "'0.0875"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.31%  '
$ws.Range('D34').Value = "'5.28"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('D35').Value = "'0.123"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.37%  '
$ws.Range('E36').Value = '  +14.44%  '
$ws.Range('D37').Value = "'4.49"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('E38').Value = '  -3.90%  '
$ws.Range('D39').Value = "'13.31"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.02%  '
$ws.Range('D40').Value = "'2.86"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +18.31%  '
$ws.Range('D41').Value = "'2.15"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.01%  '
$ws.Range('D42').Value = "'0.205"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('D43').Value = "'61.87"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.11%  '
$ws.Range('D44').Value = "'5.38"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('D45').Value = "'104.34"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('D47').Value = "'0.473"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('D48').Value = "'0.0985"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').Value = "'1.13"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.13%  '
$ws.Range('E50').Value = '  +1.46%  '
$ws.Range('D51').Value = "'1.48"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +24.40%  '
